# "Add files via upload" — a new date column (09-ago) was appended to the
# daily tracker after column AU (08-ago). Populate column AV on Sheet1 with
# the header label plus each row's value, then leave the cursor where the
# author left it (AW7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in row 1 — a new shared string "09-ago" (same pattern as the
# other date headers in C1..AU1).
$ws.Range("AV1").Value = "09-ago"

# New data column values, row by row.
$ws.Range("AV2").Value  = 0
$ws.Range("AV3").Value  = 13.963963764522161
$ws.Range("AV4").Value  = 18.562568061245496
$ws.Range("AV5").Value  = 17.904324546162204
$ws.Range("AV6").Value  = 0
$ws.Range("AV7").Value  = 15.84464979576129
$ws.Range("AV8").Value  = 15.175125690339041
$ws.Range("AV9").Value  = 10.415863499507099
$ws.Range("AV10").Value = 15.281259224339818
$ws.Range("AV11").Value = 12.780891376323334
$ws.Range("AV12").Value = 0
$ws.Range("AV13").Value = 10.243537395381232
$ws.Range("AV14").Value = 0
$ws.Range("AV15").Value = 0
$ws.Range("AV16").Value = 10.71097779190605
$ws.Range("AV17").Value = 0
$ws.Range("AV18").Value = 0

# Match the author's final selection (cell cursor moved from AW5 to AW7).
$ws.Range("AW7").Select() | Out-Null
